$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New departure rows appended to the "Main Data" table (rows 279-283),
# continuing the existing NUMBER/DATE/TIME/FLIGHT/TO/SHORT/AIRLINE/MODEL/
# AIRCFAT ID/STATUS/DIFFERENCE layout.

$rows = @(
    @{ Row=279; A=278; B="Saturday, Jan 14"; C="6:20 PM"; D="LO3886"; E="Warsaw";     F="(WAW)"; G="LOT (Star Alliance Livery) "; H="E170"; I="(SP-LDK)"; J="6:07 PM";  L="0 hours, -13 minutes" },
    @{ Row=280; A=279; B="Saturday, Jan 14"; C="7:00 PM"; D="FR7100"; E="Oslo";       F="(OSL)"; G="Ryanair ";                    H="B738"; I="(SP-RSN)"; J="7:10 PM";  L="0 hours, 10 minutes" },
    @{ Row=281; A=280; B="Saturday, Jan 14"; C="8:20 PM"; D="3Z7108"; E="Marsa Alam"; F="(RMF)"; G="Smartwings ";                 H="B38M"; I="(OK-SWC)"; J="8:25 PM";  L="0 hours, 5 minutes" },
    @{ Row=282; A=281; B="Saturday, Jan 14"; C="8:44 PM"; D="LPR42";  E="Warsaw";     F="(WAW)"; G="Polish Medical Air Rescue ";  H="LJ75"; I="(SP-MXS)"; J="8:40 PM";  L="0 hours, -4 minutes" },
    @{ Row=283; A=282; B="Saturday, Jan 14"; C="9:35 PM"; D="FR6391"; E="London";     F="(STN)"; G="Ryanair ";                    H="B738"; I="(EI-GSH)"; J="10:02 PM"; L="0 hours, 27 minutes" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
    $ws.Range("L$row").Value = $r.L
}
